$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3005
$ws1.Range("F11").Value = 132

# Sheet "全部类型" (fourth sheet) mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3005
$ws4.Range("F11").Value = 132
